$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing rows 2-5 with new values
$ws.Range("A2:J6").Value = $null

$data = @(
    @(1, 0, 3, 2, 7, 2, 4, 45, 5, "train_dim2_1"),
    @(2, 1, 2, 2, 7, 1, 5, 56, 5, "train_dim2_1"),
    @(3, 1, 4, 6, 5, 5, 1, 12, 5, "train_dim2_1"),
    @(4, 0, 1, 3, 4, 3, 3, 34, 5, "train_dim2_1"),
    @(5, 3, 4, 7, 6, 4, 2, 23, 5, "train_dim2_1")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
